$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.762.47"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "2.225.39"
$ws.Range("E3").Value = "  -5.24%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "291.74"
$ws.Range("E5").Value = "  -6.27%  "
$ws.Range("D6").Value = "84.32"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("D7").Value = "0.512"
$ws.Range("E7").Value = "  -3.10%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -3.38%  "
$ws.Range("D10").Value = "0.0794"
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("D11").Value = "29.90"
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("D12").Value = "47.99"
$ws.Range("E12").Value = "  -8.60%  "
$ws.Range("E13").Value = "  -2.31%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "6.33"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.565.43"
$ws.Range("E15").Value = "  -5.27%  "
$ws.Range("D16").Value = "14.07"
$ws.Range("E16").Value = "  -5.00%  "
$ws.Range("D17").Value = "2.208.69"
$ws.Range("E17").Value = "  -6.52%  "
$ws.Range("D18").Value = "0.720"
$ws.Range("E18").Value = "  -5.49%  "
$ws.Range("D19").Value = "39.658.50"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("D20").Value = "0.0₃0887"
$ws.Range("E20").Value = "  -1.71%  "
$ws.Range("E21").Value = "  -5.47%  "
$ws.Range("D22").Value = "65.16"
$ws.Range("E22").Value = "  -4.62%  "
$ws.Range("D23").Value = "10.42"
$ws.Range("E23").Value = "  -2.15%  "
$ws.Range("D24").Value = "231.84"
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "2.39"
$ws.Range("E26").Value = "  -6.30%  "
$ws.Range("D27").Value = "1.83"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "22.85"
$ws.Range("E28").Value = "  -3.55%  "
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("D31").Value = "154.28"
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").Value = "32.60"
$ws.Range("E32").Value = "  -6.62%  "
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").Value = "4.82"
$ws.Range("E34").Value = "  -5.92%  "
$ws.Range("E35").Value = "  -5.25%  "
$ws.Range("E36").Value = "  -2.48%  "
$ws.Range("D37").Value = "16.06"
$ws.Range("E37").Value = "  +2.41%  "
$ws.Range("E38").Value = "  -2.07%  "
$ws.Range("D39").Value = "0.0978"
$ws.Range("E40").Value = "  -5.21%  "
$ws.Range("D41").Value = "1.65"
$ws.Range("D42").Value = "3.71"
$ws.Range("E42").Value = "  -3.92%  "
$ws.Range("D43").Value = "1.947.05"
$ws.Range("E43").Value = "  -1.36%  "
$ws.Range("E44").Value = "  -3.42%  "
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("D46").Value = "9.28"
$ws.Range("E46").Value = "  -2.63%  "
$ws.Range("D47").Value = "16.12"
$ws.Range("E47").Value = "  -8.19%  "
$ws.Range("E48").Value = "  -4.46%  "
$ws.Range("D49").Value = "2.437.54"
$ws.Range("E49").Value = "  -5.06%  "
$ws.Range("D50").Value = "70.62"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").Value = "88.88"
$ws.Range("E51").Value = "  -4.76%  "
